$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.10123517574713919
$ws.Range("B3").Value = -0.00006589318244967414
$ws.Range("C3").Value = 0.0005674880895330798
$ws.Range("D3").Value = 0.9795135246293617
$ws.Range("E3").Value = 0.03370482344440404
$ws.Range("F3").Value = -0.001178152894107712
$ws.Range("G3").Value = 0.001046366529208363
$ws.Range("H3").Value = 0.1011692825646895
$ws.Range("B4").Value = 0.00291162297504885
$ws.Range("C4").Value = 0.0007469998063396009
$ws.Range("D4").Value = 5.78359035317301373
$ws.Range("E4").Value = 0.0609781558597781
$ws.Range("F4").Value = 0.001447525661577126
$ws.Range("G4").Value = 0.004375720288520576
$ws.Range("H4").Value = 0.104146798722188
$ws.Range("B5").Value = 0.006762449917766807
$ws.Range("C5").Value = 0.004721457749015145
$ws.Range("D5").Value = 5.93025575077966405
$ws.Range("E5").Value = 0.14678088976802661
$ws.Range("F5").Value = -0.0024914650115488
$ws.Range("G5").Value = 0.01601636484708241
$ws.Range("H5").Value = 0.10799762566490601
$ws.Range("B6").Value = 0.02143571216739928
$ws.Range("C6").Value = 0.003003103530011584
$ws.Range("D6").Value = 8.97147740901482038
$ws.Range("E6").Value = 0.03025110050235567
$ws.Range("F6").Value = 0.01554971729684603
$ws.Range("G6").Value = 0.02732170703795254
$ws.Range("H6").Value = 0.1226708879145384
$ws.Range("B7").Value = 0.03151197749777614
$ws.Range("C7").Value = 0.007496828470972115
$ws.Range("D7").Value = 9.4773690483778541
$ws.Range("E7").Value = 0.07222834080136985
$ws.Range("F7").Value = 0.01681841592586954
$ws.Range("G7").Value = 0.04620553906968274
$ws.Range("H7").Value = 0.13274715324491529
$ws.Range("B8").Value = 0.01566130156890224
$ws.Range("C8").Value = 0.00486468554622695
$ws.Range("D8").Value = 10.31196551574564069
$ws.Range("E8").Value = 0.1669509256071664
$ws.Range("F8").Value = 0.006126667363332324
$ws.Range("G8").Value = 0.02519593577447217
$ws.Range("H8").Value = 0.1168964773160414
$ws.Range("B9").Value = 0.01888128793421891
$ws.Range("C9").Value = 0.004313933376101708
$ws.Range("D9").Value = 9.50168468968510993
$ws.Range("E9").Value = 0.06484465993604548
$ws.Range("F9").Value = 0.0104261099769643
$ws.Range("G9").Value = 0.02733646589147354
$ws.Range("H9").Value = 0.1201164636813581
$ws.Range("B10").Value = -0.10123517574713919
$ws.Range("C10").Value = 0.0004569170895974157
$ws.Range("D10").Value = -231.27102807707859711
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -0.1021307196290461
$ws.Range("G10").Value = -0.1003396318652322
$ws.Range("B11").Value = -0.04664907511981412
$ws.Range("C11").Value = 0.0005022436580067734
$ws.Range("D11").Value = -94.69447320746908758
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = -0.04763345772349965
$ws.Range("G11").Value = -0.04566469251612856
$ws.Range("H11").Value = 0.05458610062732504
$ws.Range("B12").Value = -0.03564490283051065
$ws.Range("C12").Value = 0.0004930296048573858
$ws.Range("D12").Value = -73.40546459695268311
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = -0.03661122617224972
$ws.Range("G12").Value = -0.03467857948877159
$ws.Range("H12").Value = 0.0655902729166285
$ws.Range("B13").Value = -0.03121984147822765
$ws.Range("C13").Value = 0.0004846639115231738
$ws.Range("D13").Value = -66.34369168636906977
$ws.Range("E13").Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000007149166547085229
$ws.Range("F13").Value = -0.03216976830032803
$ws.Range("G13").Value = -0.03026991465612728
$ws.Range("H13").Value = 0.0700153342689115
$ws.Range("B14").Value = -0.02484604011828673
$ws.Range("C14").Value = 0.0004719036128126717
$ws.Range("D14").Value = -54.69510651557462211
$ws.Range("E14").Value = 0.0000000001387789208973135
$ws.Range("F14").Value = -0.02577095712519416
$ws.Range("G14").Value = -0.0239211231113793
$ws.Range("H14").Value = 0.07638913562885241
$ws.Range("B15").Value = -0.02161352500604998
$ws.Range("C15").Value = 0.0004636887179979733
$ws.Range("D15").Value = -47.81952074182017753
$ws.Range("E15").Value = 0.0000000000000000000000000000000000000000000006176359900274828
$ws.Range("F15").Value = -0.02252234107827477
$ws.Range("G15").Value = -0.0207047089338252
$ws.Range("H15").Value = 0.07962165074108918
$ws.Range("B16").Value = -0.01889730664974518
$ws.Range("C16").Value = 0.0004542128935453514
$ws.Range("D16").Value = -42.66986255050220223
$ws.Range("E16").Value = 0.00000000000000000000000000578117023347495
$ws.Range("F16").Value = -0.01978755038488428
$ws.Range("G16").Value = -0.01800706291460608
$ws.Range("H16").Value = 0.08233786909739398
$ws.Range("B17").Value = -0.01842734644764261
$ws.Range("C17").Value = 0.0004616591031851987
$ws.Range("D17").Value = -41.34445837815888325
$ws.Range("E17").Value = 0.0000000000000000000000000000002092169898259992
$ws.Range("F17").Value = -0.01933218452491676
$ws.Range("G17").Value = -0.01752250837036846
$ws.Range("H17").Value = 0.08280782929949654
$ws.Range("B18").Value = -0.01576213634729972
$ws.Range("C18").Value = 0.0004671509827339276
$ws.Range("D18").Value = -35.3336181370218867
$ws.Range("E18").Value = 0.0000000002735245380952394
$ws.Range("F18").Value = -0.01667773834423398
$ws.Range("G18").Value = -0.01484653435036545
$ws.Range("H18").Value = 0.08547303939983944
$ws.Range("B19").Value = -0.01289638177448963
$ws.Range("C19").Value = 0.0004642816788500543
$ws.Range("D19").Value = -29.91969099351907957
$ws.Range("E19").Value = 0.005533691752497484
$ws.Range("F19").Value = -0.01380636002867284
$ws.Range("G19").Value = -0.01198640352030644
$ws.Range("H19").Value = 0.08833879397264952
$ws.Range("B20").Value = -0.009556246252433917
$ws.Range("C20").Value = 0.0004732218380064631
$ws.Range("D20").Value = -21.53070293156033088
$ws.Range("E20").Value = 0.04313611207010481
$ws.Range("F20").Value = -0.01048374694876245
$ws.Range("G20").Value = -0.008628745556105375
$ws.Range("H20").Value = 0.09167892949470524
$ws.Range("B21").Value = -0.007483090728984607
$ws.Range("C21").Value = 0.0004678575007426358
$ws.Range("D21").Value = -17.14187201933362914
$ws.Range("E21").Value = 0.004006807459472007
$ws.Range("F21").Value = -0.008400077466540227
$ws.Range("G21").Value = -0.006566103991428987
$ws.Range("H21").Value = 0.09375208501815455
$ws.Range("B22").Value = -0.006643100040246571
$ws.Range("C22").Value = 0.0004600386074716878
$ws.Range("D22").Value = -15.29792069801091969
$ws.Range("E22").Value = 0.08041343124995798
$ws.Range("F22").Value = -0.007544761984005525
$ws.Range("G22").Value = -0.005741438096487619
$ws.Range("H22").Value = 0.09459207570689258
$ws.Range("B23").Value = -0.004521917638630506
$ws.Range("C23").Value = 0.0004544467219598246
$ws.Range("D23").Value = -10.30063886520706973
$ws.Range("E23").Value = 0.07351305597777417
$ws.Range("F23").Value = -0.005412619638724751
$ws.Range("G23").Value = -0.003631215638536261
$ws.Range("H23").Value = 0.09671325810850864
$ws.Range("B24").Value = -0.004219403220280883
$ws.Range("C24").Value = 0.0004524574446649575
$ws.Range("D24").Value = -9.39966856183902522
$ws.Range("E24").Value = 0.09135848328437413
$ws.Range("F24").Value = -0.005106206310932456
$ws.Range("G24").Value = -0.003332600129629311
$ws.Range("H24").Value = 0.09701577252685827
$ws.Range("B25").Value = -0.002246222561597904
$ws.Range("C25").Value = 0.0004403535532082086
$ws.Range("D25").Value = -5.02704411111781369
$ws.Range("E25").Value = 0.08629822316105711
$ws.Range("F25").Value = -0.003109302392396154
$ws.Range("G25").Value = -0.001383142730799654
$ws.Range("H25").Value = 0.09898895318554125
$ws.Range("B26").Value = 0.01546207237526988
$ws.Range("C26").Value = 0.001301714124591056
$ws.Range("D26").Value = 17.26932551464404142
$ws.Range("E26").Value = 0.06454088899839001
$ws.Range("F26").Value = 0.01291075199247867
$ws.Range("G26").Value = 0.01801339275806111
$ws.Range("H26").Value = 0.116697248122409
